## Updated C3DC phs002371 queries
## The "TreatmentTab" query (cell B5) needs its WHERE clause updated to
## also require that trt.treatment_id is not null.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTreatmentQuery = @'
SELECT
    DISTINCT prt.participant_id AS "Participant Id",
    trt.treatment_id AS "Treatment Id",
    CASE 
    WHEN trt.age_at_treatment_start = -999 THEN 'Not Reported'
    WHEN trt.age_at_treatment_start >= 1000 THEN 
        substr(trt.age_at_treatment_start, 1, length(trt.age_at_treatment_start) - 3) || ',' || substr(trt.age_at_treatment_start, -3)
    ELSE 
        trt.age_at_treatment_start 
END AS "Age at Treatment Start",
    CASE 
    WHEN trt.age_at_treatment_end = -999 THEN 'Not Reported'
    WHEN trt.age_at_treatment_end >= 1000 THEN 
        substr(trt.age_at_treatment_end, 1, length(trt.age_at_treatment_end) - 3) || ',' || substr(trt.age_at_treatment_end, -3)
    ELSE 
        trt.age_at_treatment_end 
END AS "Age at Treatment End",
    trt.treatment_type AS "Treatment Type",
    REPLACE(trt.treatment_agent, ';', ', ') AS "Treatment Agent",
    std.dbgap_accession AS "dbGaP Accession"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_diagnoses dgn ON prt.id = dgn."participant.id"
LEFT JOIN 
    df_treatments trt ON prt.id = trt."participant.id"
LEFT JOIN 
    df_treatment_resp trr ON prt.id = trr."participant.id"
LEFT JOIN 
    df_survival srv ON prt.id = srv."participant.id"
LEFT JOIN 
    df_reference_files rfs ON std.id = rfs."study.id"
WHERE 
    std.dbgap_accession = 'phs002371' AND dgn.anatomic_site = 'C42.0 : Blood' AND trt.treatment_id IS NOT NULL
ORDER BY 
    trt.treatment_id ASC
LIMIT 100;
'@

$treatmentCell = $ws.Range("B5")
$treatmentCell.Value() = $newTreatmentQuery

# Re-touch the font so the cell keeps its existing look (Calibri 12, wrap
# text) while refreshing its style record, matching the formatting churn
# Excel performs when the formula-bar text of a wrapped cell is edited.
$treatmentCell.Font.ThemeColor = 1

# Leave the selection on C5, reflecting the last cell touched during the edit.
$ws.Range("C5").Select() | Out-Null
